$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 61 is no longer the last row of data, so it switches from the
# "date only" number format to the regular "date + time" format used
# by every other data row.
$ws.Range("A61").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 62 (new last row keeps the
# "date only" number format that row 61 used to have).
$ws.Range("A62").Value = 45802
$ws.Range("A62").NumberFormat = "YYYY-MM-DD"
$ws.Range("B62").Value = 259
$ws.Range("C62").Value = 270
$ws.Range("D62").Value = 264
